$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 260
$ws1.Range("F4").Value = 2666
$ws1.Range("F5").Value = 49
$ws1.Range("F6").Value = 568

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 260
$ws4.Range("F6").Value = 2666
$ws4.Range("F7").Value = 49
$ws4.Range("F8").Value = 568
